$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in / correct the "Absent" column (H) values to form the consolidated report
$ws.Range("H4").Value = 1
$ws.Range("H7").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 0
